$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix A9: "CFP10 SR10" -> "CFP2 SR10"
$ws.Range("A9").Value = "CFP2 SR10"

# Add the two new FDR rows (description first, then their part numbers)
$ws.Range("A34").Value = "FDR"
$ws.Range("A35").Value = "FDR Gen3"
$ws.Range("C34").Value = "FCBN414QB1;FCBG414QB1"
$ws.Range("C35").Value = "FCBN414QD3;FCCN414QD3"

# Trim the Quadwire Gen3 part-number list in C19
$ws.Range("C19").Value = "FCCG410QD3;FCBG410QD3;FCBN410QD3;FCCN410QD3"

# Append the remaining new rows 36-40 (description + part numbers)
$ws.Range("A36").Value = "FDR transiver"
$ws.Range("C36").Value = "FTL414QB2;FTL414QL2"

$ws.Range("A37").Value = "SNAP12"
$ws.Range("C37").Value = "FTXD02SL1"

$ws.Range("A38").Value = "Octopus"
$ws.Range("C38").Value = "FCBR510QE2;FCBN510QE2"

$ws.Range("A39").Value = "QSFPSR4 FET Gen2"
$ws.Range("C39").Value = "FTL410QT2"

$ws.Range("A40").Value = "QSFPSR4 FET Gen3"
$ws.Range("C40").Value = "FTL410QT3"

# Scroll/selection state to match the final view
$ws.Range("C40").Select()
$excel.ActiveWindow.ScrollRow = 19
